# Refresh the cryptos list: latest prices + 1h volume deltas scraped
# from coinranking.com, including two rank swaps (rows 15/16 and 28/29).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.118.55"
$ws.Range("E2").Value = "  +2.08%  "

$ws.Range("D3").Value = "2.533.63"
$ws.Range("E3").Value = "  +0.40%  "

$ws.Range("E4").Value = "  +0.00%  "

$prevStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.70"
$ws.Range("D5").Style = $prevStyle
$ws.Range("E5").Value = "  +1.50%  "

$prevStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.05"
$ws.Range("D6").Style = $prevStyle
$ws.Range("E6").Value = "  +3.25%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D9").Value = "2.532.49"
$ws.Range("E9").Value = "  +0.42%  "

$prevStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("D10").Style = $prevStyle
$ws.Range("E10").Value = "  +1.81%  "

$ws.Range("E11").Value = "  +2.67%  "

$ws.Range("E12").Value = "  -0.30%  "

$ws.Range("E13").Value = "  -1.13%  "

$prevStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.83"
$ws.Range("D14").Style = $prevStyle
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.993.39"
$ws.Range("E15").Value = "  +0.26%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$prevStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000180"
$ws.Range("D16").Style = $prevStyle
$ws.Range("E16").Value = "  +2.17%  "

$ws.Range("D17").Value = "67.995.84"
$ws.Range("E17").Value = "  +2.13%  "

$ws.Range("D18").Value = "2.533.64"
$ws.Range("E18").Value = "  +0.37%  "

$prevStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.99"
$ws.Range("D19").Style = $prevStyle
$ws.Range("E19").Value = "  +1.64%  "

$prevStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.55"
$ws.Range("D20").Style = $prevStyle
$ws.Range("E20").Value = "  +1.88%  "

$prevStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "367.15"
$ws.Range("D21").Style = $prevStyle
$ws.Range("E21").Value = "  +5.50%  "

$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("E23").Value = "  +0.98%  "

$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("E25").Value = "  -2.35%  "

$prevStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.90"
$ws.Range("D26").Style = $prevStyle
$ws.Range("E26").Value = "  +0.95%  "

$prevStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.19"
$ws.Range("D27").Style = $prevStyle
$ws.Range("E27").Value = "  +2.44%  "

$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$prevStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = $prevStyle
$ws.Range("E28").Value = "  +0.44%  "

$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.659.92"
$ws.Range("E29").Value = "  +0.97%  "

$ws.Range("D30").Value = "0.0₃0997"
$ws.Range("E30").Value = "  +1.62%  "

$prevStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "541.40"
$ws.Range("D31").Style = $prevStyle
$ws.Range("E31").Value = "  +2.59%  "

$ws.Range("E32").Value = "  +1.60%  "

$ws.Range("E33").Value = "  +1.17%  "

$ws.Range("E34").Value = "  +1.83%  "

$ws.Range("E35").Value = "  -1.34%  "

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("E37").Value = "  -0.33%  "

$prevStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.23"
$ws.Range("D38").Style = $prevStyle
$ws.Range("E38").Value = "  +0.09%  "

$prevStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.86"
$ws.Range("D39").Style = $prevStyle
$ws.Range("E39").Value = "  +1.16%  "

$ws.Range("E40").Value = "  +1.66%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("E42").Value = "  +0.56%  "

$prevStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.20"
$ws.Range("D43").Style = $prevStyle
$ws.Range("E43").Value = "  +2.03%  "

$ws.Range("E44").Value = "  +1.49%  "

$ws.Range("E45").Value = "  -0.10%  "

$prevStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.31"
$ws.Range("D46").Style = $prevStyle
$ws.Range("E46").Value = "  -1.46%  "

$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("E48").Value = "  +0.92%  "

$ws.Range("E49").Value = "  +2.73%  "

$ws.Range("E50").Value = "  -1.07%  "

$ws.Range("E51").Value = "  -0.26%  "

